# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the batch of files that were just handed off, and sets their "Priority" to "ht".

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 13)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-05 22:29:51"
}

# zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-09-05 22:29:45"
    $wsZhCn.Range("E$r").Value = "ht"
}

# de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-09-05 22:29:51"
    $wsDeDe.Range("E$r").Value = "ht"
}
